$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 0.42
$ws.Range("E2").Value = 0.5
$ws.Range("F2").Value = 0.98
$ws.Range("G2").Value = 0.31

# Row 3 updates
$ws.Range("B3").Value = 0.35
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 0.42
$ws.Range("E3").ClearContents()
$ws.Range("F3").ClearContents()
$ws.Range("G3").ClearContents()

# Row 4 updates
$ws.Range("B4").Value = 0.35
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 0.42
